$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 1; $r -le 8; $r++) {
    $ws.Cells.Item($r, 2).Value = $r
    $ws.Cells.Item($r, 3).Value = $r
}

$ws.Range("C1:C8").Select()
